# Update the "Comment" column text for a few of the bear-quiz answers, and
# flip the "Correct" flag for the Polar row from Y to N, matching the
# re-uploaded version of this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Wrong! Totally overrated."
$ws.Range("B3").Value = "N"
$ws.Range("C3").Value = "Wrong!  They eat penguins."
$ws.Range("C4").Value = "Wrong!  It's the same as a grizzly.  Isn't it?  Or not?"
$ws.Range("C6").Value = "It's not even a bear, for pete's sake."

# Leave the active selection on C6, matching the saved state of the file.
$ws.Range("C6").Select()
